# BP-359 Bank excel statemenst upload
# Re-key the reconciliation sheet: TRN ref numbers (col B) become text
# values read straight off the bank statement instead of numbers, and
# apply the formatting that goes with the new "statement upload" layout
# (text format for the reference columns, 2-decimal format for amounts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: TRN_REF_NO becomes text, new values from the statement ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "199999"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "288888"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "388888"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "488888"

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "588888"

# --- Column A (bank codes) and C (ref no) also get the text format ---
$ws.Range("A2:A6").NumberFormat = "@"
$ws.Range("C2:C6").NumberFormat = "@"

# --- Column E (txn amount) gets a 2-decimal number format ---
$ws.Range("E2:E6").NumberFormat = "0.00"

# --- Print as portrait ---
$ws.PageSetup.Orientation = 1

# --- Move the active selection ---
$ws.Range("C3").Select() | Out-Null
